$wb = $excel.ActiveWorkbook

# ===== Sheet: yearly =====
$ws1 = $wb.Worksheets.Item("yearly")

# Header row: copy header style from G1 into H1:K1 first
$ws1.Range("G1").Copy() | Out-Null
$ws1.Range("H1:K1").PasteSpecial(-4122) | Out-Null

$ws1.Cells.Item(1,1).Value = "year"
$ws1.Cells.Item(1,2).Value = "ElasticNet_RMSE"
$ws1.Cells.Item(1,3).Value = "ElasticNet_MAE"
$ws1.Cells.Item(1,4).Value = "KNeighborsRegressor_RMSE"
$ws1.Cells.Item(1,5).Value = "KNeighborsRegressor_MAE"
$ws1.Cells.Item(1,6).Value = "RandomForest_RMSE"
$ws1.Cells.Item(1,7).Value = "RandomForest_MAE"
$ws1.Cells.Item(1,8).Value = "XGBoost_RMSE"
$ws1.Cells.Item(1,9).Value = "XGBoost_MAE"
$ws1.Cells.Item(1,10).Value = "MLP_RMSE"
$ws1.Cells.Item(1,11).Value = "MLP_MAE"

# Data rows
$ws1.Cells.Item(2,1).Value = 2018
$ws1.Cells.Item(2,2).Value = 156.5236700530149
$ws1.Cells.Item(2,3).Value = 73.71748597211828
$ws1.Cells.Item(2,4).Value = 169.2274168628197
$ws1.Cells.Item(2,5).Value = 66.33229747298925
$ws1.Cells.Item(2,6).Value = 154.2681548989847
$ws1.Cells.Item(2,7).Value = 60.45944334548784
$ws1.Cells.Item(2,8).Value = 170.28000263353
$ws1.Cells.Item(2,9).Value = 70.11495971679688
$ws1.Cells.Item(2,10).Value = 192.5173633065867
$ws1.Cells.Item(2,11).Value = 68.64278558515392

$ws1.Cells.Item(3,1).Value = 2019
$ws1.Cells.Item(3,2).Value = 116.0262056294069
$ws1.Cells.Item(3,3).Value = 56.54756130685386
$ws1.Cells.Item(3,4).Value = 104.3126050507923
$ws1.Cells.Item(3,5).Value = 43.99988472166719
$ws1.Cells.Item(3,6).Value = 105.4346927801412
$ws1.Cells.Item(3,7).Value = 42.20581722707435
$ws1.Cells.Item(3,8).Value = 133.7159565319338
$ws1.Cells.Item(3,9).Value = 51.25881958007812
$ws1.Cells.Item(3,10).Value = 113.8506388186275
$ws1.Cells.Item(3,11).Value = 49.06485462590436

# ===== Sheet: monthly =====
$ws2 = $wb.Worksheets.Item("monthly")

# Header row: copy header style from H1 into I1:L1 first
$ws2.Range("H1").Copy() | Out-Null
$ws2.Range("I1:L1").PasteSpecial(-4122) | Out-Null

$ws2.Cells.Item(1,1).Value = "year"
$ws2.Cells.Item(1,2).Value = "month"
$ws2.Cells.Item(1,3).Value = "ElasticNet_RMSE"
$ws2.Cells.Item(1,4).Value = "ElasticNet_MAE"
$ws2.Cells.Item(1,5).Value = "KNeighborsRegressor_RMSE"
$ws2.Cells.Item(1,6).Value = "KNeighborsRegressor_MAE"
$ws2.Cells.Item(1,7).Value = "RandomForest_RMSE"
$ws2.Cells.Item(1,8).Value = "RandomForest_MAE"
$ws2.Cells.Item(1,9).Value = "XGBoost_RMSE"
$ws2.Cells.Item(1,10).Value = "XGBoost_MAE"
$ws2.Cells.Item(1,11).Value = "MLP_RMSE"
$ws2.Cells.Item(1,12).Value = "MLP_MAE"

# Data rows
$ws2.Cells.Item(2,1).Value = 2018
$ws2.Cells.Item(2,2).Value = 4
$ws2.Cells.Item(2,3).Value = 161.3191211485808
$ws2.Cells.Item(2,4).Value = 81.19825776809301
$ws2.Cells.Item(2,5).Value = 173.07922243973
$ws2.Cells.Item(2,6).Value = 73.81867780513075
$ws2.Cells.Item(2,7).Value = 155.778144304089
$ws2.Cells.Item(2,8).Value = 60.14513126199558
$ws2.Cells.Item(2,9).Value = 167.5928031159154
$ws2.Cells.Item(2,10).Value = 70.61150360107422
$ws2.Cells.Item(2,11).Value = 210.477261747838
$ws2.Cells.Item(2,12).Value = 75.75588175822841

$ws2.Cells.Item(3,2).Value = 5
$ws2.Cells.Item(3,3).Value = 173.4026094816599
$ws2.Cells.Item(3,4).Value = 81.13077554938047
$ws2.Cells.Item(3,5).Value = 192.6322064650257
$ws2.Cells.Item(3,6).Value = 72.7178380922722
$ws2.Cells.Item(3,7).Value = 167.9337783814861
$ws2.Cells.Item(3,8).Value = 59.22335320979651
$ws2.Cells.Item(3,9).Value = 184.2310348957254
$ws2.Cells.Item(3,10).Value = 70.17586517333984
$ws2.Cells.Item(3,11).Value = 202.7635361631993
$ws2.Cells.Item(3,12).Value = 76.59191660908748

$ws2.Cells.Item(4,2).Value = 6
$ws2.Cells.Item(4,3).Value = 151.7971728731756
$ws2.Cells.Item(4,4).Value = 74.8734505682569
$ws2.Cells.Item(4,5).Value = 175.921188740304
$ws2.Cells.Item(4,6).Value = 58.74223596825557
$ws2.Cells.Item(4,7).Value = 147.2995759864605
$ws2.Cells.Item(4,8).Value = 60.62198652287365
$ws2.Cells.Item(4,9).Value = 163.7527015521423
$ws2.Cells.Item(4,10).Value = 67.17802429199219
$ws2.Cells.Item(4,11).Value = 193.3533356839787
$ws2.Cells.Item(4,12).Value = 73.60765126952012

$ws2.Cells.Item(5,2).Value = 7
$ws2.Cells.Item(5,3).Value = 172.7749225385216
$ws2.Cells.Item(5,4).Value = 72.99231065197414
$ws2.Cells.Item(5,5).Value = 173.5888976674047
$ws2.Cells.Item(5,6).Value = 63.63691719118808
$ws2.Cells.Item(5,7).Value = 168.761412352992
$ws2.Cells.Item(5,8).Value = 64.2884743517012
$ws2.Cells.Item(5,9).Value = 180.8008767435877
$ws2.Cells.Item(5,10).Value = 69.19409942626953
$ws2.Cells.Item(5,11).Value = 175.4396583451101
$ws2.Cells.Item(5,12).Value = 65.73653637429634

$ws2.Cells.Item(6,2).Value = 8
$ws2.Cells.Item(6,3).Value = 189.5309205992445
$ws2.Cells.Item(6,4).Value = 79.46960028095748
$ws2.Cells.Item(6,5).Value = 190.2438563676492
$ws2.Cells.Item(6,6).Value = 74.93244875017871
$ws2.Cells.Item(6,7).Value = 185.7253069313445
$ws2.Cells.Item(6,8).Value = 65.64733621001963
$ws2.Cells.Item(6,9).Value = 200.2761179908878
$ws2.Cells.Item(6,10).Value = 74.98692321777344
$ws2.Cells.Item(6,11).Value = 152.0103710881688
$ws2.Cells.Item(6,12).Value = 63.621972119213

$ws2.Cells.Item(7,2).Value = 9
$ws2.Cells.Item(7,3).Value = 183.8967028575
$ws2.Cells.Item(7,4).Value = 85.49214438970562
$ws2.Cells.Item(7,5).Value = 207.0844366902592
$ws2.Cells.Item(7,6).Value = 85.96245335532063
$ws2.Cells.Item(7,7).Value = 159.9322551462637
$ws2.Cells.Item(7,8).Value = 69.52937656381332
$ws2.Cells.Item(7,9).Value = 189.0513219836085
$ws2.Cells.Item(7,10).Value = 81.79102325439453
$ws2.Cells.Item(7,11).Value = 207.9255974384463
$ws2.Cells.Item(7,12).Value = 78.32595362482428

$ws2.Cells.Item(8,2).Value = 10
$ws2.Cells.Item(8,3).Value = 137.5177804102361
$ws2.Cells.Item(8,4).Value = 69.90864997680514
$ws2.Cells.Item(8,5).Value = 150.5139965622938
$ws2.Cells.Item(8,6).Value = 64.24480342767728
$ws2.Cells.Item(8,7).Value = 180.0839217584433
$ws2.Cells.Item(8,8).Value = 65.87632854629354
$ws2.Cells.Item(8,9).Value = 171.3080307698679
$ws2.Cells.Item(8,10).Value = 75.25382995605469
$ws2.Cells.Item(8,11).Value = 230.6171320669188
$ws2.Cells.Item(8,12).Value = 71.58136852908561

$ws2.Cells.Item(9,2).Value = 11
$ws2.Cells.Item(9,3).Value = 115.2175207623537
$ws2.Cells.Item(9,4).Value = 64.02767366744028
$ws2.Cells.Item(9,5).Value = 124.8634506714973
$ws2.Cells.Item(9,6).Value = 56.61465219084417
$ws2.Cells.Item(9,7).Value = 111.418748592552
$ws2.Cells.Item(9,8).Value = 54.66055904505739
$ws2.Cells.Item(9,9).Value = 133.4353702991639
$ws2.Cells.Item(9,10).Value = 66.43912506103516
$ws2.Cells.Item(9,11).Value = 217.3889575040643
$ws2.Cells.Item(9,12).Value = 65.70094522539367

$ws2.Cells.Item(10,2).Value = 12
$ws2.Cells.Item(10,3).Value = 101.0675219850619
$ws2.Cells.Item(10,4).Value = 54.82827170557587
$ws2.Cells.Item(10,5).Value = 113.1747613538537
$ws2.Cells.Item(10,6).Value = 47.21063857390361
$ws2.Cells.Item(10,7).Value = 83.45419609665466
$ws2.Cells.Item(10,8).Value = 44.49420189840944
$ws2.Cells.Item(10,9).Value = 130.0792201530475
$ws2.Cells.Item(10,10).Value = 55.81435775756836
$ws2.Cells.Item(10,11).Value = 113.431670147239
$ws2.Cells.Item(10,12).Value = 46.86189922446385

$ws2.Cells.Item(11,1).Value = 2019
$ws2.Cells.Item(11,2).Value = 1
$ws2.Cells.Item(11,3).Value = 94.90714842927271
$ws2.Cells.Item(11,4).Value = 57.17775423234387
$ws2.Cells.Item(11,5).Value = 93.55475462935391
$ws2.Cells.Item(11,6).Value = 45.84110991998527
$ws2.Cells.Item(11,7).Value = 111.1594989781828
$ws2.Cells.Item(11,8).Value = 53.36473688017569
$ws2.Cells.Item(11,9).Value = 120.8340642937661
$ws2.Cells.Item(11,10).Value = 56.51005172729492
$ws2.Cells.Item(11,11).Value = 149.2974677605195
$ws2.Cells.Item(11,12).Value = 64.49988538257743

$ws2.Cells.Item(12,2).Value = 2
$ws2.Cells.Item(12,3).Value = 115.1614386917969
$ws2.Cells.Item(12,4).Value = 56.90678874500144
$ws2.Cells.Item(12,5).Value = 106.1042845181113
$ws2.Cells.Item(12,6).Value = 45.44493331715852
$ws2.Cells.Item(12,7).Value = 115.7395708278442
$ws2.Cells.Item(12,8).Value = 42.14743828218399
$ws2.Cells.Item(12,9).Value = 140.8420963126792
$ws2.Cells.Item(12,10).Value = 52.06947708129883
$ws2.Cells.Item(12,11).Value = 100.6285405313758
$ws2.Cells.Item(12,12).Value = 44.45297501954031

$ws2.Cells.Item(13,2).Value = 3
$ws2.Cells.Item(13,3).Value = 115.02086386918
$ws2.Cells.Item(13,4).Value = 60.82994716638951
$ws2.Cells.Item(13,5).Value = 100.5091048137546
$ws2.Cells.Item(13,6).Value = 47.60684918849587
$ws2.Cells.Item(13,7).Value = 113.6376969496651
$ws2.Cells.Item(13,8).Value = 45.61263593944486
$ws2.Cells.Item(13,9).Value = 141.4802741330041
$ws2.Cells.Item(13,10).Value = 55.35153579711914
$ws2.Cells.Item(13,11).Value = 102.6804891141668
$ws2.Cells.Item(13,12).Value = 51.60529877460304

$ws2.Cells.Item(14,2).Value = 4
$ws2.Cells.Item(14,3).Value = 137.3372865546844
$ws2.Cells.Item(14,4).Value = 60.33880971893009
$ws2.Cells.Item(14,5).Value = 114.1560353009774
$ws2.Cells.Item(14,6).Value = 48.37883126470116
$ws2.Cells.Item(14,7).Value = 112.6169553102496
$ws2.Cells.Item(14,8).Value = 44.19515546148872
$ws2.Cells.Item(14,9).Value = 148.1067565242214
$ws2.Cells.Item(14,10).Value = 53.45566177368164
$ws2.Cells.Item(14,11).Value = 119.5119213278615
$ws2.Cells.Item(14,12).Value = 50.00508633651857

$ws2.Cells.Item(15,2).Value = 5
$ws2.Cells.Item(15,3).Value = 141.8838933305355
$ws2.Cells.Item(15,4).Value = 61.60928092231516
$ws2.Cells.Item(15,5).Value = 109.8864493395742
$ws2.Cells.Item(15,6).Value = 46.22391842594758
$ws2.Cells.Item(15,7).Value = 110.8690343362119
$ws2.Cells.Item(15,8).Value = 47.21504510449822
$ws2.Cells.Item(15,9).Value = 150.5878130962961
$ws2.Cells.Item(15,10).Value = 56.10960388183594
$ws2.Cells.Item(15,11).Value = 111.1501196871345
$ws2.Cells.Item(15,12).Value = 48.97425342335239

$ws2.Cells.Item(16,2).Value = 6
$ws2.Cells.Item(16,3).Value = 125.6666502822947
$ws2.Cells.Item(16,4).Value = 58.13117600920312
$ws2.Cells.Item(16,5).Value = 98.23977085798951
$ws2.Cells.Item(16,6).Value = 39.3775633894932
$ws2.Cells.Item(16,7).Value = 94.57080990833721
$ws2.Cells.Item(16,8).Value = 36.61422699604585
$ws2.Cells.Item(16,9).Value = 131.2516145734025
$ws2.Cells.Item(16,10).Value = 49.20993804931641
$ws2.Cells.Item(16,11).Value = 117.439039480862
$ws2.Cells.Item(16,12).Value = 45.67584579989976

$ws2.Cells.Item(17,2).Value = 7
$ws2.Cells.Item(17,3).Value = 108.7597354114627
$ws2.Cells.Item(17,4).Value = 51.13001185988653
$ws2.Cells.Item(17,5).Value = 108.2934599669307
$ws2.Cells.Item(17,6).Value = 39.25323893325756
$ws2.Cells.Item(17,7).Value = 97.27004814309707
$ws2.Cells.Item(17,8).Value = 36.63343632080069
$ws2.Cells.Item(17,9).Value = 127.6351726785665
$ws2.Cells.Item(17,10).Value = 47.02345275878906
$ws2.Cells.Item(17,11).Value = 106.6236865570407
$ws2.Cells.Item(17,12).Value = 46.16594943260831

$ws2.Cells.Item(18,2).Value = 8
$ws2.Cells.Item(18,3).Value = 92.38655882158278
$ws2.Cells.Item(18,4).Value = 49.98999853907212
$ws2.Cells.Item(18,5).Value = 107.6871928562987
$ws2.Cells.Item(18,6).Value = 40.9139161600098
$ws2.Cells.Item(18,7).Value = 87.75323217210331
$ws2.Cells.Item(18,8).Value = 35.65833416327474
$ws2.Cells.Item(18,9).Value = 112.4236589938368
$ws2.Cells.Item(18,10).Value = 44.40501403808594
$ws2.Cells.Item(18,11).Value = 101.1473642069431
$ws2.Cells.Item(18,12).Value = 43.13970224570934

$ws2.Cells.Item(19,2).Value = 9
$ws2.Cells.Item(19,3).Value = 99.20127096738939
$ws2.Cells.Item(19,4).Value = 50.99381972075057
$ws2.Cells.Item(19,5).Value = 99.20364086902475
$ws2.Cells.Item(19,6).Value = 41.73775472727318
$ws2.Cells.Item(19,7).Value = 97.00237524968978
$ws2.Cells.Item(19,8).Value = 35.87146162855679
$ws2.Cells.Item(19,9).Value = 121.4365993275915
$ws2.Cells.Item(19,10).Value = 45.0150260925293
$ws2.Cells.Item(19,11).Value = 103.8566115276451
$ws2.Cells.Item(19,12).Value = 45.04487124417369

# ===== Sheet: per_split =====
$ws3 = $wb.Worksheets.Item("per_split")

# Header row: copy header style from G1 into H1:K1 first
$ws3.Range("G1").Copy() | Out-Null
$ws3.Range("H1:K1").PasteSpecial(-4122) | Out-Null

$ws3.Cells.Item(1,1).Value = "split_index"
$ws3.Cells.Item(1,2).Value = "ElasticNet_RMSE"
$ws3.Cells.Item(1,3).Value = "ElasticNet_MAE"
$ws3.Cells.Item(1,4).Value = "KNeighborsRegressor_RMSE"
$ws3.Cells.Item(1,5).Value = "KNeighborsRegressor_MAE"
$ws3.Cells.Item(1,6).Value = "RandomForest_RMSE"
$ws3.Cells.Item(1,7).Value = "RandomForest_MAE"
$ws3.Cells.Item(1,8).Value = "XGBoost_RMSE"
$ws3.Cells.Item(1,9).Value = "XGBoost_MAE"
$ws3.Cells.Item(1,10).Value = "MLP_RMSE"
$ws3.Cells.Item(1,11).Value = "MLP_MAE"

# Data rows
$ws3.Cells.Item(2,1).Value = 1
$ws3.Cells.Item(2,2).Value = 162.3474127638753
$ws3.Cells.Item(2,3).Value = 79.02811554608327
$ws3.Cells.Item(2,4).Value = 180.7527455157689
$ws3.Cells.Item(2,5).Value = 68.33307209904014
$ws3.Cells.Item(2,6).Value = 157.1719164991475
$ws3.Cells.Item(2,7).Value = 60.00047239812633
$ws3.Cells.Item(2,8).Value = 172.0529873868803
$ws3.Cells.Item(2,9).Value = 69.30071258544922
$ws3.Cells.Item(2,10).Value = 202.2184698830858
$ws3.Cells.Item(2,11).Value = 75.30447808282881

$ws3.Cells.Item(3,1).Value = 2
$ws3.Cells.Item(3,2).Value = 182.1891755289138
$ws3.Cells.Item(3,3).Value = 79.29700667058214
$ws3.Cells.Item(3,4).Value = 190.7393186232258
$ws3.Cells.Item(3,5).Value = 74.8062400136719
$ws3.Cells.Item(3,6).Value = 171.8346110655409
$ws3.Cells.Item(3,7).Value = 66.47873541941746
$ws3.Cells.Item(3,8).Value = 190.20645649859
$ws3.Cells.Item(3,9).Value = 75.30238342285156
$ws3.Cells.Item(3,10).Value = 179.8424449797392
$ws3.Cells.Item(3,11).Value = 69.20157719385298

$ws3.Cells.Item(4,1).Value = 3
$ws3.Cells.Item(4,2).Value = 118.9259715508949
$ws3.Cells.Item(4,3).Value = 62.94275889184711
$ws3.Cells.Item(4,4).Value = 130.4923816617192
$ws3.Cells.Item(4,5).Value = 56.0462535108739
$ws3.Cells.Item(4,6).Value = 131.508095682494
$ws3.Cells.Item(4,7).Value = 55.03728831016112
$ws3.Cells.Item(4,8).Value = 146.1761164870992
$ws3.Cells.Item(4,9).Value = 65.86177825927734
$ws3.Cells.Item(4,10).Value = 194.5178825513076
$ws3.Cells.Item(4,11).Value = 61.42034445002332

$ws3.Cells.Item(5,1).Value = 4
$ws3.Cells.Item(5,2).Value = 108.7410069431822
$ws3.Cells.Item(5,3).Value = 58.3029173749476
$ws3.Cells.Item(5,4).Value = 100.1662077959782
$ws3.Cells.Item(5,5).Value = 46.29707119486381
$ws3.Cells.Item(5,6).Value = 113.5198514072021
$ws3.Cells.Item(5,7).Value = 47.06190045730217
$ws3.Cells.Item(5,8).Value = 134.6875797563569
$ws3.Cells.Item(5,9).Value = 54.6505241394043
$ws3.Cells.Item(5,10).Value = 119.7656306811699
$ws3.Cells.Item(5,11).Value = 53.55505906053833

$ws3.Cells.Item(6,1).Value = 5
$ws3.Cells.Item(6,2).Value = 135.1605710810613
$ws3.Cells.Item(6,3).Value = 60.03142262381062
$ws3.Cells.Item(6,4).Value = 107.6697984954359
$ws3.Cells.Item(6,5).Value = 44.67855529247279
$ws3.Cells.Item(6,6).Value = 106.3665856536819
$ws3.Cells.Item(6,7).Value = 42.69162570011937
$ws3.Cells.Item(6,8).Value = 143.607946115457
$ws3.Cells.Item(6,9).Value = 52.93473815917969
$ws3.Cells.Item(6,10).Value = 116.0904269924753
$ws3.Cells.Item(6,11).Value = 48.22727080502125

$ws3.Cells.Item(7,1).Value = 6
$ws3.Cells.Item(7,2).Value = 100.4536087094986
$ws3.Cells.Item(7,3).Value = 50.70470700172046
$ws3.Cells.Item(7,4).Value = 105.2680936001423
$ws3.Cells.Item(7,5).Value = 40.60276016685766
$ws3.Cells.Item(7,6).Value = 94.10160832683046
$ws3.Cells.Item(7,7).Value = 36.06354814765468
$ws3.Cells.Item(7,8).Value = 120.7258726542285
$ws3.Cells.Item(7,9).Value = 45.50517272949219
$ws3.Cells.Item(7,10).Value = 103.930545811875
$ws3.Cells.Item(7,11).Value = 44.7944195519255

# Row 8: AVERAGE formulas extended through column K
$ws3.Cells.Item(8,2).Formula = "=AVERAGE(B2:B7)"
$ws3.Cells.Item(8,3).Formula = "=AVERAGE(C2:C7)"
$ws3.Cells.Item(8,4).Formula = "=AVERAGE(D2:D7)"
$ws3.Cells.Item(8,5).Formula = "=AVERAGE(E2:E7)"
$ws3.Cells.Item(8,6).Formula = "=AVERAGE(F2:F7)"
$ws3.Cells.Item(8,7).Formula = "=AVERAGE(G2:G7)"
$ws3.Cells.Item(8,8).Formula = "=AVERAGE(H2:H7)"
$ws3.Cells.Item(8,9).Formula = "=AVERAGE(I2:I7)"
$ws3.Cells.Item(8,10).Formula = "=AVERAGE(J2:J7)"
$ws3.Cells.Item(8,11).Formula = "=AVERAGE(K2:K7)"

# Selection update recorded in the diff
$ws3.Activate() | Out-Null
$ws3.Range("G13").Select() | Out-Null

Write-Output "done"